$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 165; existing rows 165-169 shift down to 166-170.
$ws.Rows("165").Insert()

# Populate the newly inserted row 165 with the new weekly data point.
$ws.Cells(165, 1).Value = 9
$ws.Cells(165, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells(165, 3).Value = "Metropolitana"
$ws.Cells(165, 4).Value = 44509
$ws.Cells(165, 5).Value = 13
$ws.Cells(165, 6).Value = 100112026
$ws.Cells(165, 7).Value = "Haba"
$ws.Cells(165, 8).Value = "Sin especificar"
$ws.Cells(165, 9).Value = "Primera"
$ws.Cells(165, 10).Value = 52
$ws.Cells(165, 11).Value = 7000
$ws.Cells(165, 12).Value = 8000
$ws.Cells(165, 13).Value = 7500
$ws.Cells(165, 14).Value = '$/saco 25 kilos'
$ws.Cells(165, 15).Value = "Región Metropolitana"
$ws.Cells(165, 16).Value = 300
$ws.Cells(165, 17).Value = 25
$ws.Cells(165, 18).Value = "Hortaliza"
